$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.63
$ws.Range("I3").Value = 2.55
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 9.5
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 3.4
$ws.Range("W3").Value = 8
$ws.Range("Z3").Value = 26
$ws.Range("AA3").Value = 21
$ws.Range("AK3").Value = 26
# Row 4
$ws.Range("G4").Value = 2.9
$ws.Range("I4").Value = 2.4
$ws.Range("J4").Value = 3.5
$ws.Range("L4").Value = 3
$ws.Range("W4").Value = 9.5
$ws.Range("X4").Value = 15
$ws.Range("AJ4").Value = 9.5
$ws.Range("AL4").Value = 19
$ws.Range("AN4").Value = 5
$ws.Range("BA4").Value = 51
# Row 5
$ws.Range("G5").Value = 2.55
$ws.Range("I5").Value = 2.7
$ws.Range("J5").Value = 3.1
$ws.Range("L5").Value = 3.25
$ws.Range("W5").Value = 9.5
$ws.Range("X5").Value = 13
$ws.Range("AA5").Value = 19
$ws.Range("AB5").Value = 26
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 15
$ws.Range("AL5").Value = 21
$ws.Range("AM5").Value = 29
$ws.Range("AO5").Value = 13
$ws.Range("AP5").Value = 21
$ws.Range("AQ5").Value = 41
$ws.Range("AX5").Value = 15
$ws.Range("AY5").Value = 23
$ws.Range("AZ5").Value = 51
# Row 6
$ws.Range("G6").Value = 1.6
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 5.25
$ws.Range("K6").Value = 2.3
$ws.Range("Q6").Value = 1.8
$ws.Range("R6").Value = 2
$ws.Range("X6").Value = 8
$ws.Range("Z6").Value = 12
$ws.Range("AA6").Value = 13
$ws.Range("AO6").Value = 8
$ws.Range("AP6").Value = 19
# Row 7
$ws.Range("G7").Value = 2
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 2.63
$ws.Range("L7").Value = 4
$ws.Range("W7").Value = 8
$ws.Range("X7").Value = 10
$ws.Range("AC7").Value = 12
$ws.Range("AO7").Value = 11
$ws.Range("AU7").Value = 7.5
$ws.Range("AX7").Value = 19
$ws.Range("BB7").Value = 151
$ws.Range("BC7").Value = 151
# Row 8
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
# Row 9
$ws.Range("J9").Value = 1.4
$ws.Range("K9").Value = 3.5
$ws.Range("L9").Value = 12
$ws.Range("S9").Value = 1.17
$ws.Range("T9").Value = 5
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.73
$ws.Range("AN9").Value = 3.5
$ws.Range("AP9").Value = 15
$ws.Range("AQ9").Value = 9
$ws.Range("AR9").Value = 26
$ws.Range("AS9").Value = 81
$ws.Range("AT9").Value = 5
$ws.Range("AV9").Value = 51
$ws.Range("AY9").Value = 41
$ws.Range("AZ9").Value = 301
$ws.Range("BA9").Value = 201
$ws.Range("BB9").Value = 301
# Row 10
$ws.Range("G10").Value = 3.2
$ws.Range("H10").Value = 2.88
$ws.Range("I10").Value = 2.45
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = 1.91
$ws.Range("L10").Value = 3.25
$ws.Range("W10").Value = 7.5
$ws.Range("X10").Value = 15
$ws.Range("Y10").Value = 13
$ws.Range("Z10").Value = 34
$ws.Range("AA10").Value = 34
$ws.Range("AB10").Value = 41
$ws.Range("AC10").Value = 6
$ws.Range("AD10").Value = 5.5
$ws.Range("AE10").Value = 17
$ws.Range("AH10").Value = 6.5
$ws.Range("AI10").Value = 10
$ws.Range("AJ10").Value = 11
$ws.Range("AK10").Value = 23
$ws.Range("AL10").Value = 23
$ws.Range("AN10").Value = 5
$ws.Range("AO10").Value = 19
$ws.Range("AQ10").Value = 67
$ws.Range("AR10").Value = 101
$ws.Range("AS10").Value = 301
$ws.Range("AW10").Value = 4.33
$ws.Range("AX10").Value = 15
$ws.Range("AZ10").Value = 51
$ws.Range("BB10").Value = 301
# Row 11
$ws.Range("G11").Value = 1.91
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 2.5
$ws.Range("L11").Value = 4
$ws.Range("O11").Value = 1.22
$ws.Range("P11").Value = 4
$ws.Range("Q11").Value = 1.75
$ws.Range("R11").Value = 2.05
$ws.Range("W11").Value = 9
$ws.Range("X11").Value = 10
$ws.Range("Z11").Value = 17
$ws.Range("AA11").Value = 15
$ws.Range("AC11").Value = 12
$ws.Range("AD11").Value = 6.5
$ws.Range("AE11").Value = 12
$ws.Range("AH11").Value = 13
$ws.Range("AI11").Value = 21
$ws.Range("AJ11").Value = 13
$ws.Range("AL11").Value = 29
$ws.Range("AO11").Value = 10
$ws.Range("AQ11").Value = 34
$ws.Range("AR11").Value = 51
$ws.Range("AS11").Value = 126
$ws.Range("AX11").Value = 19
